# Update the timestamp embedded in the test e-mail addresses on the
# "UsuariosRegistro" sheet from 20251109_013943 to 20251109_020650.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("UsuariosRegistro")
$ws2 = $wb.Worksheets.Item("LoginData")

$ws1.Range("C2").Value = "juan.perez+20251109_020650@test.com"
$ws1.Range("C3").Value = "maria.gonzalez+20251109_020650@test.com"
$ws1.Range("C4").Value = "carlos.rodriguez+20251109_020650@test.com"
$ws1.Range("C5").Value = "ana.martinez+20251109_020650@test.com"
$ws1.Range("C6").Value = "luis.garcia+20251109_020650@test.com"

$ws2.Range("A2").Value = "juan.perez+20251109_020650@test.com"
$ws2.Range("A3").Value = "maria.gonzalez+20251109_020650@test.com"
